# Update V2G_scenarios.xlsx:
#  - Change the "V2G mandate" scenario ramp: 2030->2028 (0.12), add 2031 (0.5)
#    and 2035 (0.88), keep 2050 but change value to 0.9.
#  - Add a new "Early" scenario with the same year/value ramp as the
#    updated "V2G mandate" scenario (2020/0, 2024/0.12, 2027/0.5, 2031/0.88, 2050/0.9).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing "V2G mandate" rows (rows 12-13) and insert two more ---

# Row 12: year 2030 -> 2028, value 1 -> 0.12
$ws.Cells.Item(12, 3).Value = 2028
$ws.Cells.Item(12, 4).Value = 0.12

# Row 13 used to be (2050, 1) for "V2G mandate" - turn it into (2031, 0.5)
$ws.Cells.Item(13, 3).Value = 2031
$ws.Cells.Item(13, 4).Value = 0.5

# Insert two new rows after row 13 for "V2G mandate" (2035/0.88) and keep a
# final (2050/0.9) row, pushing "No V2G" data down.
$ws.Rows.Item(14).Insert()
$ws.Rows.Item(14).Insert()

$ws.Cells.Item(14, 1).Value = "V2G mandate"
$ws.Cells.Item(14, 2).Value = "BEV"
$ws.Cells.Item(14, 3).Value = 2035
$ws.Cells.Item(14, 4).Value = 0.88

$ws.Cells.Item(15, 1).Value = "V2G mandate"
$ws.Cells.Item(15, 2).Value = "BEV"
$ws.Cells.Item(15, 3).Value = 2050
$ws.Cells.Item(15, 4).Value = 0.9

# --- "No V2G" scenario now lives at rows 16-18 (shifted down by 2) ---
$ws.Cells.Item(16, 1).Value = "No V2G"
$ws.Cells.Item(16, 2).Value = "BEV"
$ws.Cells.Item(16, 3).Value = 2020
$ws.Cells.Item(16, 4).Value = 0

$ws.Cells.Item(17, 1).Value = "No V2G"
$ws.Cells.Item(17, 2).Value = "BEV"
$ws.Cells.Item(17, 3).Value = 2030
$ws.Cells.Item(17, 4).Value = 0

$ws.Cells.Item(18, 1).Value = "No V2G"
$ws.Cells.Item(18, 2).Value = "BEV"
$ws.Cells.Item(18, 3).Value = 2050
$ws.Cells.Item(18, 4).Value = 0

# --- Append the new "Early" scenario (rows 19-23) ---
$ws.Cells.Item(19, 1).Value = "Early"
$ws.Cells.Item(19, 2).Value = "BEV"
$ws.Cells.Item(19, 3).Value = 2020
$ws.Cells.Item(19, 4).Value = 0

$ws.Cells.Item(20, 1).Value = "Early"
$ws.Cells.Item(20, 2).Value = "BEV"
$ws.Cells.Item(20, 3).Value = 2024
$ws.Cells.Item(20, 4).Value = 0.12

$ws.Cells.Item(21, 1).Value = "Early"
$ws.Cells.Item(21, 2).Value = "BEV"
$ws.Cells.Item(21, 3).Value = 2027
$ws.Cells.Item(21, 4).Value = 0.5

$ws.Cells.Item(22, 1).Value = "Early"
$ws.Cells.Item(22, 2).Value = "BEV"
$ws.Cells.Item(22, 3).Value = 2031
$ws.Cells.Item(22, 4).Value = 0.88

$ws.Cells.Item(23, 1).Value = "Early"
$ws.Cells.Item(23, 2).Value = "BEV"
$ws.Cells.Item(23, 3).Value = 2050
$ws.Cells.Item(23, 4).Value = 0.9

# Move the active selection to match the source workbook's saved view.
$ws.Range("C20").Select()
